# Generate Report for Handback
# Updates the handoff/handback timestamps for the second tracked file
# (c33b6286-7556-4d32-9fc9-7304f58620f1) across the Overview, zh-cn and
# de-de worksheets, reflecting a newly generated handback report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for the c33b6286 file (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-21 22:56:47"

# --- zh-cn sheet: Correspond Handoff / Handback datetimes for the c33b6286 file (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-21 22:56:43"
$wsZhCn.Range("K3").Value = "2016-08-21 22:56:59"

# --- de-de sheet: Correspond Handoff / Handback datetimes for the c33b6286 file (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-21 22:56:47"
$wsDeDe.Range("K3").Value = "2016-08-21 22:57:10"
